# Upgraded to Live Scraping and Universal App ID support
# Refreshes the trend report workbook with the latest scraped mention
# counts (Trend Analysis), and recomputes the derived Summary /
# Trending Topics / Top Topics sheets from that new data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Trend Analysis" — raw daily mention counts per topic
# (columns B..AF = 2025-11-26 .. 2025-12-26), one row per topic.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Trend Analysis")

$trend = @{
    2 = @(3,0,0,3,4,1,1,3,3,0,4,1,3,3,1,5,2,2,1,4,2,0,6,4,0,2,6,2,1,6,2)
    3 = @(1,2,7,8,4,2,4,5,5,4,2,4,5,2,2,1,2,5,5,3,3,3,5,6,7,4,4,1,5,5,2)
    4 = @(2,2,3,2,5,1,2,0,0,1,7,7,2,5,1,3,2,5,3,2,2,3,1,3,5,3,5,3,4,4,3)
    5 = @(4,2,2,1,1,1,1,3,3,4,2,1,0,4,0,2,5,2,2,4,3,4,2,1,2,4,4,4,2,2,1)
    6 = @(0,3,1,1,7,2,2,0,2,1,2,4,0,5,1,2,3,2,3,3,2,3,2,1,2,4,4,0,4,2,1)
    7 = @(3,1,1,4,1,4,4,4,3,2,3,5,7,5,1,2,2,5,3,2,1,2,2,3,6,1,6,2,2,4,2)
    8 = @(2,3,3,1,2,1,4,3,2,0,4,4,1,1,2,5,4,3,4,2,1,2,4,5,4,4,0,3,1,1,3)
}

foreach ($r in $trend.Keys) {
    $vals = $trend[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        # column B is index 2
        $ws1.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
}

# ---------------------------------------------------------------
# Sheet 2: "Summary" — headline metrics derived from the new data
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B3").Value = 598          # Total Mentions
$ws2.Range("B4").Value = 85.43        # Average Daily Mentions
$ws2.Range("B5").Value = "Customer support unresponsive"  # Most Frequent Topic
$ws2.Range("B7").Value = 1            # Topics with Significant Growth

# ---------------------------------------------------------------
# Sheet 3: "Trending Topics" — week-over-week growth, re-ranked
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Trending Topics")

$trendingTopics = @(
    @{ Row=2; Topic="Delivery issue";                 Last=3.86; Prev=2.71; Rate=0.42;  Pct="42.1%";  Total=91  },
    @{ Row=3; Topic="Payment issue";                  Last=3.29; Prev=2.57; Rate=0.28;  Pct="27.8%";  Total=93  },
    @{ Row=4; Topic="Delivery partner rude";           Last=2.71; Prev=2.57; Rate=0.06;  Pct="5.6%";   Total=73  },
    @{ Row=5; Topic="Food stale";                     Last=2.43; Prev=2.29; Rate=0.06;  Pct="6.2%";   Total=69  },
    @{ Row=6; Topic="App crashing";                   Last=2.71; Prev=2.71; Rate=0;     Pct="0.0%";   Total=75  },
    @{ Row=7; Topic="Customer support unresponsive";  Last=4;    Prev=4.29; Rate=-0.07; Pct="-6.7%";  Total=118 },
    @{ Row=8; Topic="Wrong order delivered";           Last=2.29; Prev=3;    Rate=-0.24; Pct="-23.8%"; Total=79  }
)

foreach ($entry in $trendingTopics) {
    $r = $entry.Row
    $ws3.Cells.Item($r, 1).Value = $entry.Topic
    $ws3.Cells.Item($r, 2).Value = $entry.Last
    $ws3.Cells.Item($r, 3).Value = $entry.Prev
    $ws3.Cells.Item($r, 4).Value = $entry.Rate

    # "Growth %" is stored as a literal text string (e.g. "42.1%"), not a
    # numeric percentage — force text interpretation so Excel doesn't
    # auto-convert it to a fraction, then drop back to the default style
    # so we don't leave a stray number-format override on the cell.
    $pctCell = $ws3.Cells.Item($r, 5)
    $pctCell.NumberFormat = "@"
    $pctCell.Value = $entry.Pct
    $pctCell.Style = "Normal"

    $ws3.Cells.Item($r, 6).Value = $entry.Total
}

# ---------------------------------------------------------------
# Sheet 4: "Top Topics" — re-ranked by total mentions, descending
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Top Topics")

$topTopics = @(
    @{ Row=2; Topic="Customer support unresponsive"; Total=118; Avg=3.81 },
    @{ Row=3; Topic="Payment issue";                 Total=93;  Avg=3    },
    @{ Row=4; Topic="Delivery issue";                Total=91;  Avg=2.94 },
    @{ Row=5; Topic="Wrong order delivered";          Total=79;  Avg=2.55 },
    @{ Row=6; Topic="App crashing";                  Total=75;  Avg=2.42 },
    @{ Row=7; Topic="Delivery partner rude";          Total=73;  Avg=2.35 },
    @{ Row=8; Topic="Food stale";                    Total=69;  Avg=2.23 }
)

foreach ($entry in $topTopics) {
    $r = $entry.Row
    $ws4.Cells.Item($r, 1).Value = $entry.Topic
    $ws4.Cells.Item($r, 2).Value = $entry.Total
    $ws4.Cells.Item($r, 3).Value = $entry.Avg
}
